$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) "Joint Master & PhD" line: add "since" before the year 2024
#    "...School of Computing<TAB> 2024" -> "...School of Computing<TAB>since 2024"
#    (search only the text AFTER the tab stop so the <w:tab/> run itself is
#     left untouched by the replace)
# ---------------------------------------------------------------
$tab = [char]9
$eduPara = $d.Paragraphs(6)
$tabRange = $eduPara.Range.Duplicate
$tabRange.Find.Execute($tab, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$afterTab = $d.Range($tabRange.End, $eduPara.Range.End)
$afterTab.Find.Execute(" 2024", $true, $false, $false, $false, $false, $true, 1, $false, "since 2024", 2) | Out-Null

# ---------------------------------------------------------------
# 2) "Research Interest" line
#    "Applied AI " -> "Human-centered computing "
#    "Mobile, IoT and Wearable Computing " -> "Mobile Computing · Applied AI"
# ---------------------------------------------------------------
$r2 = $d.Content
$r2.Find.Execute("Applied AI ", $true, $false, $false, $false, $false, $true, 1, $false, "Human-centered computing ", 2) | Out-Null

$r3 = $d.Content
$r3.Find.Execute("Mobile, IoT and Wearable Computing ", $true, $false, $false, $false, $false, $true, 1, $false, "Mobile Computing · Applied AI", 2) | Out-Null

# ---------------------------------------------------------------
# 3) "Bachelor ..." line: drop the "Graduated in " label before the year
#    "...4.00<TAB>Graduated in 2022" -> "...4.00<TAB>2022"
# ---------------------------------------------------------------
$r4 = $d.Content
$r4.Find.Execute("Graduated in ", $true, $false, $false, $false, $false, $true, 1, $false, "", 2) | Out-Null

# ---------------------------------------------------------------
# 4) VTM internship paragraph: rewrite the description
# ---------------------------------------------------------------
$oldVtm = "Visual Token Matching (VTM) is a general-purpose few-shot learner for arbitrary visual dense prediction tasks, as proposed by a lab mate in an outstanding paper in ICLR’23. However, VTM cannot handle temporal information, which hinders its performance in video domains. In my internship, I enhanced VTM’s generalizability by incorporating time attention into its framework. Empirical results show that the method surpasses the baseline VTM when a very limited support set is available. Specifically, the method achieves 8.89% and 4.37% higher accuracy than the baseline in 1-shot and 2-shot scenarios, respectively, on the DAVIS2016 video segmentation dataset. "
$newVtm = "Visual Token Matching (ICLR’23) is a general-purpose few-shot learner for arbitrary visual dense prediction tasks. To address its limitation in modeling temporal information, I integrated temporal attention into its framework. On the DAVIS2016 benchmark, the improved model outperformed the baseline by 8.89% in 1-shot and 4.37% in 2-shot settings. "
$r5 = $d.Content
$r5.Find.Execute($oldVtm, $true, $false, $false, $false, $false, $true, 1, $false, $newVtm, 2) | Out-Null

# ---------------------------------------------------------------
# 5) Undergraduate research paragraph: rewrite the description
# ---------------------------------------------------------------
$dash = [char]8211
$oldTrack = "We created a novel inference architecture that leverages re-identification features for data association in visual object tracking for long-term videos. Our tracker provisionally matched the state-of-the-art performance within the scope of person tracking in the Visual Object Tracking ${dash} Long Term 2021 benchmark."
$newTrack = "We developed a novel inference architecture using re-identification features for data association in long-term visual object tracking. Our tracker achieved competitive performance on person tracking in the VOT-LT 2021 benchmark. "
$r6 = $d.Content
$r6.Find.Execute($oldTrack, $true, $false, $false, $false, $false, $true, 1, $false, $newTrack, 2) | Out-Null

Write-Host "Done"
